$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Row 1: turn the old (buggy, data-filled) header row into real column headers,
# matching the header row used on the other sheets (土地/建物), plus a
# car-specific "capacity" column in place of "area".
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

# New header cells H1:N1 need the same bold/bordered header style as B1:G1.
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Row 2: keep the existing car data (A2:G2) as-is and fill in the
# previously-missing metadata columns, mirroring 土地/建物.
$ws.Range("G2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "'2012-05-01"
$ws.Range("K2").Value = "管碧玲"
$ws.Range("L2").Value = 1374
$ws.Range("M2").Value = "tmpf0df1"
$ws.Range("N2").Value = 44

# Re-apply the plain data style to J2 so the apostrophe-forced text entry
# doesn't leave a stray "quote prefix" style behind.
$ws.Range("G2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
